$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('C2').Value = '25°'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58%'
$ws.Range('H2').Value = '35°'
$ws.Range('J2').NumberFormat = '@'
$ws.Range('J2').Value = '30%'
$ws.Range('C3').Value = '25°'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '60%'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '84%'
$ws.Range('I3').Value = '20°'
$ws.Range('J3').NumberFormat = '@'
$ws.Range('J3').Value = '30%'
$ws.Range('K3').NumberFormat = '@'
$ws.Range('K3').Value = '52%'
$ws.Range('B4').Value = '35°'
$ws.Range('C4').Value = '25°'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '62%'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '84%'
$ws.Range('H4').Value = '28°'
$ws.Range('I4').Value = '15°'
$ws.Range('J4').NumberFormat = '@'
$ws.Range('J4').Value = '62%'
$ws.Range('K4').NumberFormat = '@'
$ws.Range('K4').Value = '83%'
$ws.Range('B5').Value = '33°'
$ws.Range('C5').Value = '25°'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '67%'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '87%'
$ws.Range('J5').NumberFormat = '@'
$ws.Range('J5').Value = '70%'
$ws.Range('K5').NumberFormat = '@'
$ws.Range('K5').Value = '84%'
$ws.Range('L5').Value = '8 de 11'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '72%'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '89%'
$ws.Range('J6').NumberFormat = '@'
$ws.Range('J6').Value = '66%'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '88%'
$ws.Range('K7').NumberFormat = '@'
$ws.Range('K7').Value = '68%'
$ws.Range('C8').Value = '24°'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '76%'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '93%'
$ws.Range('I8').Value = '22°'
$ws.Range('J8').NumberFormat = '@'
$ws.Range('J8').Value = '35%'
$ws.Range('K8').NumberFormat = '@'
$ws.Range('K8').Value = '49%'
$ws.Range('B9').Value = '33°'
$ws.Range('C9').Value = '24°'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '72%'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '93%'
$ws.Range('H9').Value = '33°'
$ws.Range('J9').NumberFormat = '@'
$ws.Range('J9').Value = '37%'
$ws.Range('B10').Value = '32°'
$ws.Range('C10').Value = '24°'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '73%'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '94%'
$ws.Range('H10').Value = '31°'
$ws.Range('J10').NumberFormat = '@'
$ws.Range('J10').Value = '49%'
$ws.Range('K10').NumberFormat = '@'
$ws.Range('K10').Value = '74%'
$ws.Range('B11').Value = '32°'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '72%'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '90%'
$ws.Range('H11').Value = '31°'
$ws.Range('J11').NumberFormat = '@'
$ws.Range('J11').Value = '50%'
